$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value for each data row (rows 2-176).
# The date needs to move forward by one day (from 2023-10-04 / 45203 to 2023-10-05 / 45204).
$ws.Range("C2:C176").Value = 45204
